$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the numeric values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0

# Set the label text (goes into shared strings)
$ws.Range("B2").Value = "disconnected_elements"

# Apply formatting (bold font, thin box border, centered horizontal, top
# vertical) to B1 first, fully - this builds exactly one new style entry.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.LineStyle = 1        # xlContinuous
$b1.Borders.Weight = 2           # xlThin

# Copy that exact style onto A2 (format-only paste reuses the style index
# instead of re-deriving it property-by-property).
$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)          # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
